# Apply content/structure changes to LOT2039 worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Structural changes first (shape must match before final values are written) ---
# Row 17 loses its B/C cells (shift remaining B:C content up one row)
$ws.Range("B17:C17").Delete(-4162)  # xlShiftUp
# Former row 22 (Bibliografia long text) is now empty/duplicated -> remove it entirely
$ws.Rows.Item(22).Delete()

# --- Content fixes ---

# Row 10
$ws.Range("A10").Value = 'Objetivos:'
$ws.Range("B10").Value = '3380737 - Flávio Teixeira da Silva'
$ws.Range("C10").Value = '3380737 - Flávio Teixeira da Silva'
$ws.Rows.Item(10).RowHeight = 60

# Row 13
$ws.Range("A13").Value = 'Programa resumido:'
$ws.Range("B13").Value = 'Semestral'
$ws.Range("C13").Value = 'Semestral'
$ws.Rows.Item(13).RowHeight = 60

# Row 14
$ws.Range("A14").Value = 'Short syllabus:'
$ws.Range("B14").Value = 'Structure and ultrastructure of lignocellulosic materials, cellulose, hemicelluloses, other polyoses, lignin, extractives and bark composition. Reactions in acid and alkaline media. Composition and chemical analysis of woods.'
$ws.Range("C14").Value = 'Structure and ultrastructure of lignocellulosic materials, cellulose, hemicelluloses, other polyoses, lignin, extractives and bark composition. Reactions in acid and alkaline media. Composition and chemical analysis of woods.'
$ws.Rows.Item(14).RowHeight = 60

# Row 15
$ws.Range("A15").Value = 'Programa:'
$ws.Range("B15").Value = '01/01/2018'
$ws.Range("C15").Value = '01/01/2018'
$ws.Rows.Item(15).RowHeight = 120

# Row 16
$ws.Range("A16").Value = 'Syllabus:'
$ws.Range("B16").Value = '1.Structure and ultrastructure of lignocellulosic materials: anatomic aspects. Ultrastructure of cell wall; functional elements of conductor system.2.Cellulose: occurrence; molecular properties; constitution and configuration; cellulose in solution, chain length, molar mass, hydrogen bridges; supramolecular structure; crystalline structure; fibrillar structure.3.Hemicellulose and other polyoses: nature and classification; xylan of wood of hardwoods and conifers; xylan of other plants; supramolecular structure; mannans of wood of hardwoods and conifers; other mannans; glucan; galactan and pectin.4.Lignin: significance and occurrence; cell lignification; synthesis of monomeric unities; formation of the lignin macromolecule; aspects of decomposition of cell wall; structure and constitution; models and heterogeneity; characterization and proprieties, chemical composition and molar mass; behavior at UV and infra-red; lignin-carbohydrates complexes.5.Extractives: importance; extractives of woods of hardwoods and conifers; terpenes, fat, wax, phenols, tannins, flavonoids, etc.; inorganic compounds.6.Bark composition: anatomy; chemical composition; general analysis; cellulose; polyoses, lignin, polyphenols, suberin and extractives; inorganic compounds.7.Reactions in acid medium: general aspects, reactions of polysaccharides, hydrolysis, dehydration, oxidation; reactions of lignin; solvolysis.8.Reactions in alkaline medium: reaction of polysaccharides, hydrolysis, oxidative degradation, sugar hydrogenation; lignin reactions; selective and soft oxidative degradation; hydrogenolysis.9.Chemical composition and analysis of wood: kinds of compounds; macromolecules and low molar mass substances; wood analysis, sampling and determination of inorganics and extractives; preparation of holocellulose; isolation and determination of cellulose, polyoses and lignin.'
$ws.Range("C16").Value = '1.Structure and ultrastructure of lignocellulosic materials: anatomic aspects. Ultrastructure of cell wall; functional elements of conductor system.2.Cellulose: occurrence; molecular properties; constitution and configuration; cellulose in solution, chain length, molar mass, hydrogen bridges; supramolecular structure; crystalline structure; fibrillar structure.3.Hemicellulose and other polyoses: nature and classification; xylan of wood of hardwoods and conifers; xylan of other plants; supramolecular structure; mannans of wood of hardwoods and conifers; other mannans; glucan; galactan and pectin.4.Lignin: significance and occurrence; cell lignification; synthesis of monomeric unities; formation of the lignin macromolecule; aspects of decomposition of cell wall; structure and constitution; models and heterogeneity; characterization and proprieties, chemical composition and molar mass; behavior at UV and infra-red; lignin-carbohydrates complexes.5.Extractives: importance; extractives of woods of hardwoods and conifers; terpenes, fat, wax, phenols, tannins, flavonoids, etc.; inorganic compounds.6.Bark composition: anatomy; chemical composition; general analysis; cellulose; polyoses, lignin, polyphenols, suberin and extractives; inorganic compounds.7.Reactions in acid medium: general aspects, reactions of polysaccharides, hydrolysis, dehydration, oxidation; reactions of lignin; solvolysis.8.Reactions in alkaline medium: reaction of polysaccharides, hydrolysis, oxidative degradation, sugar hydrogenation; lignin reactions; selective and soft oxidative degradation; hydrogenolysis.9.Chemical composition and analysis of wood: kinds of compounds; macromolecules and low molar mass substances; wood analysis, sampling and determination of inorganics and extractives; preparation of holocellulose; isolation and determination of cellulose, polyoses and lignin.'
$ws.Rows.Item(16).RowHeight = 120

# Row 17
$ws.Range("A17").Value = 'Avaliação:'
$ws.Rows.Item(17).AutoFit()

# Row 18
$ws.Range("A18").Value = 'Método:'
$ws.Range("B18").Value = '3380737 - Flávio Teixeira da Silva'
$ws.Range("C18").Value = '3380737 - Flávio Teixeira da Silva'
$ws.Rows.Item(18).RowHeight = 60

# Row 19
$ws.Range("A19").Value = 'Critério:'
$ws.Range("B19").Value = 'A avaliação será feita por duas provas (P1 e P2).'
$ws.Range("C19").Value = 'A avaliação será feita por duas provas (P1 e P2).'
$ws.Rows.Item(19).RowHeight = 60

# Row 20
$ws.Range("A20").Value = 'Norma de recuperação:'
$ws.Range("B20").Value = 'A nota final (NF) será calculada atribuindo-se peso um para a primeira avaliação e peso dois para a segunda avaliação do semestre.

NF=(P1 + 2xP2)/3 

Será considerado aprovado o aluno com NF>= 5,0 e 70% de freqüência no curso.'
$ws.Range("C20").Value = 'A nota final (NF) será calculada atribuindo-se peso um para a primeira avaliação e peso dois para a segunda avaliação do semestre.

NF=(P1 + 2xP2)/3 

Será considerado aprovado o aluno com NF>= 5,0 e 70% de freqüência no curso.'
$ws.Rows.Item(20).RowHeight = 60

# Row 21
$ws.Range("A21").Value = 'Bibliografia:'
$ws.Range("B21").Value = 'A recuperação será feita por meio de uma prova (PR) para alunos que tenham NF maior ou igual a 3,0 e menor do que 5,0. A nota de recuperação (NR) será calculada pela média simples entre a nota final (NF) e a prova de recuperação (PR). Será considerado aprovado o aluno com NR maior ou igual a 5,0'
$ws.Range("C21").Value = 'A recuperação será feita por meio de uma prova (PR) para alunos que tenham NF maior ou igual a 3,0 e menor do que 5,0. A nota de recuperação (NR) será calculada pela média simples entre a nota final (NF) e a prova de recuperação (PR). Será considerado aprovado o aluno com NR maior ou igual a 5,0'
$ws.Rows.Item(21).RowHeight = 120
